$d = $word.ActiveDocument

# 1) Update the two cached TIME field results from "24 de noviembre de 2024"
#    to "25 de febrero de 2025".
$d.Content.Find.Execute("24 de noviembre de 2024", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "25 de febrero de 2025", 2) | Out-Null

$d.Content.Find.Execute("24 de noviembre de 2024", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "25 de febrero de 2025", 2) | Out-Null

# 2) "el{{DIAS}}..." -> "el periodo de {{DIAS}}..."
$d.Content.Find.Execute("el{{DIAS}}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "el periodo de {{DIAS}}", 2) | Out-Null
